# LLtable.xlsx update — grammar.txt update, LLtable update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# "string" -> "str" for the existing Q1 literal cell
$ws.Range("Q1").Value = '"str"'
# New column T header
$ws.Range("T1").Value = "<expr>"

# --- Row 6: <return_expressions> gains a new T6 entry ---
$ws.Range("T6").Value = 11

# --- Row 8: <variable_definition> M8 15 -> 14 ---
$ws.Range("M8").Value = 14

# --- Row 9: <function_call> G9 16 -> 17 ---
$ws.Range("G9").Value = 17

# --- Row 10: <list_of_call_parameters> ---
$ws.Range("M10").Value = 19
$ws.Range("O10").Value = 18
$ws.Range("Q10").Value = 19
$ws.Range("S10").ClearContents()

# --- Row 11: <call_parameter> ---
$ws.Range("M11").Value = 20
$ws.Range("Q11").Value = 21

# --- Row 12: <list_of_call_parameters_n> ---
$ws.Range("O12").Value = 23
$ws.Range("P12").Value = 22
$ws.Range("S12").ClearContents()

# --- Row 13: <list_of_parameters> ---
$ws.Range("H13").Value = 25
$ws.Range("I13").Value = 25
$ws.Range("J13").Value = 25
$ws.Range("K13").Value = 25
$ws.Range("L13").Value = 25
$ws.Range("O13").Value = 24
$ws.Range("S13").ClearContents()

# --- Row 14: <parameter> ---
$ws.Range("H14").Value = 26
$ws.Range("I14").Value = 26
$ws.Range("J14").Value = 26
$ws.Range("K14").Value = 26
$ws.Range("L14").Value = 26

# --- Row 15: <list_of_parameters_n> ---
$ws.Range("O15").Value = 33
$ws.Range("P15").Value = 34
$ws.Range("S15").ClearContents()

# --- Row 16: <list_of_datatypes> ---
$ws.Range("H16").Value = 27
$ws.Range("I16").Value = 28
$ws.Range("J16").Value = 29
$ws.Range("K16").Value = 30
$ws.Range("L16").Value = 31

# --- Row 17: <variable> ---
$ws.Range("M17").Value = 32

# --- Row 18: <list_of_datatypes_ret> ---
$ws.Range("H18").Value = 36
$ws.Range("I18").Value = 36
$ws.Range("J18").Value = 36
$ws.Range("K18").Value = 36
$ws.Range("L18").Value = 36
$ws.Range("R18").Value = 35

# --- Row 19: new <var_def_expr> row ---
$ws.Range("A19").Value = "<var_def_expr>"
$ws.Range("G19").Value = 15
$ws.Range("T19").Value = 16

# --- Sheet view: zoom + selection ---
$excel.ActiveWindow.Zoom = 81
$ws.Range("R20").Select()
